# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (strikeout count replaced by K count, std/mean recalculated).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$gValues = @{
    2  = 0
    3  = 5
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    10 = 3
    11 = 0
    12 = 0
    13 = 2
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 2
    21 = 2
    22 = 0
    23 = 0
    24 = 0
    25 = 5
    26 = 0
    27 = 1
    28 = 1
    29 = 3
    30 = 0
    31 = 0
    32 = 1
    33 = 0
    34 = 2
    35 = 0
    36 = 2
    38 = 2
    39 = 2
    40 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
